# Update phone number detection logic
# - "Phone Numbers" sheet: re-detected phone numbers (col B) shuffled/reassigned per row,
#   the detected name text (col C) now also includes the phone number that was embedded
#   in the message text, and the timestamp (col D) no longer includes the time portion.
# - "Summary" sheet: same Name/Timestamp update for its single summary row.

$wb = $excel.ActiveWorkbook

$phoneSheet = $wb.Worksheets.Item("Phone Numbers")

$newPhones = @(
    "+966558927634",
    "+966531482587",
    "+966594320944",
    "+966504435170",
    "+966552914008",
    "+966590423200",
    "+966566626124",
    "+966537394446",
    "+966536276067"
)

$newName = "كلموني اخر ٢٤ ساعه, +966 55 892 7634"
$newTimestamp = "Sunday"

for ($i = 0; $i -lt $newPhones.Length; $i++) {
    $row = $i + 2
    $phoneCell = $phoneSheet.Cells.Item($row, 2)
    $phoneCell.NumberFormat = "@"
    $phoneCell.Value = $newPhones[$i]
    $phoneSheet.Cells.Item($row, 3).Value = $newName
    $phoneSheet.Cells.Item($row, 4).Value = $newTimestamp
}

$summarySheet = $wb.Worksheets.Item("Summary")
$summarySheet.Cells.Item(2, 3).Value = $newName
$summarySheet.Cells.Item(2, 4).Value = $newTimestamp
